$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 17
$ws.Range("B3").Value = 893
$ws.Range("C3").Value = 2056
$ws.Range("D3").Value = 12134

# Row 4
$ws.Range("A4").Value = 191
$ws.Range("B4").Value = 1446.6
$ws.Range("C4").Value = 5183.6000000000004
$ws.Range("D4").Value = 16791.2

# Row 6
$ws.Range("A6").Value = 52.36
$ws.Range("B6").Value = 66.02
$ws.Range("C6").Value = 71.14
$ws.Range("D6").Value = 67.709999999999994

# Row 7
$ws.Range("A7").Value = 61.84
$ws.Range("B7").Value = 66.61
$ws.Range("C7").Value = 64.42
$ws.Range("D7").Value = 59.88
